# Savesheet para diferentes bimestres
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "MCT-1A-Circuitos Elétricos 1"
$ws.Range("D3").Value = "MCT-1A-Circuitos Elétricos 1"

$ws.Range("C6").Value = "-"
$ws.Range("C7").Value = "-"

$ws.Range("B10").Value = "MEC-3A-Metrologia 2"
$ws.Range("E10").Value = "MEC-2A-Elet. Digit. Básica"
$ws.Range("F10").Value = "-"

$ws.Range("C11").Value = "-"
$ws.Range("E11").Value = "MEC-2A-Elet. Digit. Básica"
$ws.Range("F11").Value = "-"

$ws.Range("C12").Value = "-"
$ws.Range("E12").Value = "MEC-2A-Elet. Digit. Básica"

$ws.Range("C14").Value = "MEC-3A-Metrologia 2"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "MEC-2A-Elet. Digit. Básica"

$ws.Range("C15").Value = "MEC-3A-Metrologia 2"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "MCT-2A-Eletronica Analógica e de Potência"

$ws.Range("C16").Value = "MEC-3A-Metrologia 2"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "MCT-2A-Eletronica Analógica e de Potência"
